$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B" = 0.6408044419705359
    "C" = -265.1326953808737
    "D" = 0.7489112171893821
    "E" = 0.8831275631765217
    "F" = 0.8211676020406411
    "G" = 0.2132339996864685
    "H" = 157.9878643119225
    "I" = 0.1573674430258001
    "J" = 0.0845284335006489
    "K" = 0.1209479382632245
    "L" = 0.2405142646481177
    "M" = 0.4617726710043249
    "N" = 0.2163006006629874
    "O" = 0.4814312896101858
    "P" = 29.09073025240775
    "Q" = 44.93611597569436
}

foreach ($col in $values.Keys) {
    $rangeAddr = "$col" + "2:" + "$col" + "26"
    $ws.Range($rangeAddr).Value = $values[$col]
}
